# Populates the (currently blank) UX Survey document body with the full
# questionnaire text, ahead of the original sole paragraph (which carries the
# "_GoBack" bookmark and must be preserved unchanged).

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Each array entry is the *inner* content (runs / pPr) of one new top-level
# paragraph to insert, in document order. An empty string produces a bare
# empty paragraph (<w:p/>).
$paragraphs = @(
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Creative Website Solutions</w:t></w:r>',
    '<w:r><w:t>User Experience Survey</w:t></w:r>',
    '',
    '<w:r><w:t>Name:</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Date:</w:t></w:r>',
    '',
    '<w:r><w:t>What website operating system and web browser are you using?</w:t></w:r>',
    '',
    '<w:r><w:t>Did the website load quickly?</w:t></w:r>',
    '<w:r><w:t>Yes</w:t></w:r><w:r><w:tab/><w:t>No</w:t></w:r>',
    '',
    '<w:r><w:t>What is the first thing your eye is drawn to when looking at the website?</w:t></w:r>',
    '',
    '<w:r><w:t>What is your opinion of the colors of the website?</w:t></w:r>',
    '',
    '<w:r><w:t>What is your opinion of the website navigation?</w:t></w:r>',
    '',
    '<w:r><w:t>Is the website easy to read and find what you are looking for?</w:t></w:r>',
    '<w:r><w:t>Yes</w:t></w:r><w:r><w:tab/><w:t>No</w:t></w:r>',
    '',
    '<w:r><w:t>Please refer to the Test Procedures and complete the steps, note any questio</w:t></w:r><w:r><w:t>ns, concerns, or feedback in the appropriate sections below:</w:t></w:r>',
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Navigation and Design</w:t></w:r>',
    '',
    '',
    '',
    '',
    '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Registration &amp; Login</w:t></w:r>',
    '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>',
    '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>',
    '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Client Projects</w:t></w:r>',
    '',
    '',
    '',
    '',
    '',
    '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Administrator – Register, Edit, and View Users</w:t></w:r>'
)

$xml = ""
foreach ($inner in $paragraphs) {
    $xml += "<w:p $ns>$inner</w:p>"
}

# Re-append a literal copy of the document's original (only) paragraph -- the one
# carrying the _GoBack bookmark -- unchanged, including its original rsid
# attributes, so that paragraph's identity/formatting is left exactly as it was.
$tail = '<w:p ' + $ns + ' w:rsidR="00B25999" w:rsidRDefault="00B25999">' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/></w:p>'
$xml += $tail

$d = $word.ActiveDocument

# Replace the entire (currently single-paragraph) document body in one shot --
# InsertXML on a range spanning the whole body swaps in exactly this content.
# (A point insertion right before the final, bookmarked paragraph would instead
# merge the last inserted paragraph's runs into that paragraph, corrupting its
# pPr/bookmark; replacing the whole body range avoids that merge entirely.)
$full = $d.Content
$full.InsertXML($xml)

